# career_info.xlsx — add new career paths ("Culinary" senior levels + new "Law" path)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width tweaks ---
# Column C needs to get wider to fit the new, longer job titles.
$ws.Columns("C").ColumnWidth = 39.848214285714285
# Column F (old "E:F" merged width) stays the same width, but a new,
# narrow column E is introduced as a spacer between D and F.
$ws.Columns("E").ColumnWidth = 5.637274285714286
$ws.Columns("F").ColumnWidth = 39.34041428571428

# --- Existing data tweaks ---
# Row 8 & 9 no longer need the taller wrapped height now that column C is wider.
$ws.Rows(8).RowHeight = 20.05
$ws.Rows(9).RowHeight = 20.05

# CEO wage bumped from 110 to 135.
$ws.Range("D16").Value = 135

# --- New rows: copy style (borders/number format) down from the last existing row ---
$ws.Range("A22:F22").Copy()
$ws.Range("A23:F36").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row heights for the newly added rows ---
for ($r = 23; $r -le 36; $r++) {
    $ws.Rows($r).RowHeight = 20.05
}

# --- New data: two more "Culinary" levels, then the whole new "Law" path ---
$rows = @(
    @(23, "Culinary", 7,  "Restauranteur",    55,  8,   "Tuesday, Wednesday, Friday, Saturday, Sunday"),
    @(24, "Culinary", 8,  "Franchise Owner",  75,  7,   "Tuesday, Wednesday, Friday, Saturday, Sunday"),
    @(25, "Culinary", 9,  "Celebrity Chef",   100, 6.5, "Tuesday, Friday, Saturday, Sunday"),
    @(26, "Culinary", 10, "Culinary Legend",  100, 6.5, "Tuesday, Friday, Saturday, Sunday"),
    @(27, "Law", 1,  "Legal Assistant",   10,  7.5, "Monday, Tuesday, Wednesday, Thursday, Friday"),
    @(28, "Law", 2,  "Cost Draftsperson", 15,  7.5, "Monday, Tuesday, Wednesday, Thursday, Friday"),
    @(29, "Law", 3,  "Paralegal",         20,  7.5, "Monday, Tuesday, Wednesday, Thursday, Friday"),
    @(30, "Law", 4,  "Trainee Lawyer",    35,  7.5, "Monday, Tuesday, Wednesday, Thursday, Friday"),
    @(31, "Law", 5,  "Junior Associate",  55,  7.5, "Monday, Tuesday, Wednesday, Thursday, Friday"),
    @(32, "Law", 6,  "Associate",         70,  7.5, "Monday, Tuesday, Wednesday, Thursday, Friday"),
    @(33, "Law", 7,  "Senior Associate",  85,  7.5, "Monday, Tuesday, Wednesday, Thursday, Friday"),
    @(34, "Law", 8,  "Legal Director",    95,  7,   "Monday, Tuesday, Wednesday, Thursday, Friday"),
    @(35, "Law", 9,  "Partner",           110, 7,   "Monday, Tuesday, Wednesday, Thursday, Friday"),
    @(36, "Law", 10, "Senior Parter",     125, 7,   "Monday, Tuesday, Wednesday, Thursday, Friday")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
}
